$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.130.22'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.858.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.59%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.91'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.19%  '

$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4653'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.73%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2812'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.77%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06545'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.77'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07800'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.05%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -7.40%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.864.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.104'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.41%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6639'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '282.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.159.91'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.458'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.32%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.112.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007228'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.133'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.08%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.325'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.27%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.19'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.44%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.73%  '

$ws.Range("E28").Value = '  -9.90%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.340'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09553'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.418'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.470'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.101'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04637'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.84%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.098'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6985'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.24%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.707'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01835'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.328'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.44%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.507'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.88%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8526'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.917'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.14%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '103.95'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.79%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4139'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '990.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.192'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.191'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.84%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1138'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.38%  '
